$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9080969095230103
$ws.Range("B1").Value = 1.733759760856628
$ws.Range("D1").Value = 1.864870071411133
$ws.Range("E1").Value = 1.104066848754883
